$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DM_Stat (C) and P_Value (D) columns for rows 2-11
$ws.Range("C2").Value = 0.8071667876701776
$ws.Range("D2").Value = 0.4251812386364193

$ws.Range("C3").Value = 1.740132740757265
$ws.Range("D3").Value = 0.09088139420091923
$ws.Range("G3").Value = "No"

$ws.Range("C4").Value = 0.9277499347563577
$ws.Range("D4").Value = 0.3600769399336743

$ws.Range("C5").Value = 1.27481818975406
$ws.Range("D5").Value = 0.2110170255661195

$ws.Range("C6").Value = 1.574654441605272
$ws.Range("D6").Value = 0.1245967338618694

$ws.Range("C7").Value = -0.09287534776870492
$ws.Range("D7").Value = 0.926547710051268

$ws.Range("C8").Value = 0.2444921432598706
$ws.Range("D8").Value = 0.8083191949522073

$ws.Range("C9").Value = -1.055283365162098
$ws.Range("D9").Value = 0.2987371065148334

$ws.Range("C10").Value = -0.740030472408127
$ws.Range("D10").Value = 0.4643611732451853

$ws.Range("C11").Value = 0.390674774892358
$ws.Range("D11").Value = 0.698474892930788
